# Adapt slides for CS 536 in the spring.
#
# 1. The deck's cached "fixed date" footer fields (on every slide layout,
#    the slide master, and the notes master) were refreshed from
#    11/5/2024 -> 12/20/2024.
# 2. The title slide's title text was updated from the "04 -" lecture
#    number to "06,07 -" to reflect the new lecture numbering.

$p = $ppt.ActivePresentation

$oldDate = "11/5/2024"
$newDate = "12/20/2024"

function Update-DateShape($shape) {
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# Slide master's own "Date Placeholder" shape.
foreach ($shape in $p.SlideMaster.Shapes) {
    Update-DateShape $shape
}

# Every slide layout's "Date Placeholder" shape.
foreach ($layout in $p.SlideMaster.CustomLayouts) {
    foreach ($shape in $layout.Shapes) {
        Update-DateShape $shape
    }
}

# The notes master's "Date Placeholder" shape (datetimeFigureOut field).
foreach ($shape in $p.NotesMaster.Shapes) {
    Update-DateShape $shape
}

# Title slide: bump the lecture number in the title text.
$titleSlide = $p.Slides.Item(1)
$titleShape = $titleSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "06,07 – Parsing Expression Grammars + Abstract Syntax Trees"
